$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. New column F: header "overall_total"
# ------------------------------------------------------------------
$ws.Range("F1").Value = "overall_total"

# Give F1 the same base format as the existing header cells (bold font,
# thin border, centered) before the alignment tweak below is applied to
# the whole header row in a single step.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122)

# Header row: change vertical alignment from top -> center (applied to
# the whole header row A1:F1 at once so every cell lands on the same
# new style).
$ws.Range("A1:F1").VerticalAlignment = -4108

# ------------------------------------------------------------------
# 2. Data rows: add column F values + a new totals row 9, then style
#    the whole A2:F9 block (thin border + centered horiz/vert) in one
#    paste so it resolves to a single shared style.
# ------------------------------------------------------------------
$ws.Range("F9").Value = 74193

# Build the target data-row style once on a scratch cell, then copy it
# (format only) onto the whole data range in a single operation so the
# stylesheet only gains the one extra xf actually used in the diff.
$scratch = $ws.Range("H1")
$scratch.Value = 1
$scratch.Borders.LineStyle = 1
$scratch.Borders.Weight = 2
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4108

$ws.Range("H1").Copy() | Out-Null
$ws.Range("A2:F9").PasteSpecial(-4122)

$scratch.Clear()
